$d = $word.ActiveDocument

# Paragraph 1: update the date header
$d.Paragraphs(1).Range.Text = "⚡🚀 המאמר היומי של מייק 30.05.2024⚡🚀"

# Paragraph 2: replace the paper title
$null = $d.Content.Find.Execute(" Transformers Can Do Arithmetic with the Right Embeddings",
                         $true, $false, $false, $false, $false, $true, 1, $false,
                         "2BP: 2-Stage Backpropagation", 2)

# Paragraph 3: replace the review body text
$d.Paragraphs(3).Range.Text = "אנו יודעים שהמודלים העמוקים גדולים היום מדי כדי להיכנס לזיכרון ram של gpu אחד. עקב כך מחלקים את משקלי המודל בין הgpus השונים (sharding). זה פותר צוואר בקבוק אחד (זכרון) אבל כתוצאה מכך נוצר צוור בקבוק אחר בחישוב של backprop, המאמר הנסקר פיתח שיטה למקבל את חישוב הגרדיאנטים במהלך backprop ובכך מקל על צוואר הבקבוק הזה."

# Paragraph 4 ("רפו: ...") is removed entirely
$d.Paragraphs(4).Range.Delete()

# Remaining paragraph ("מאמר: ...") gets the new link
$d.Paragraphs(4).Range.Text = "מאמר: https://arxiv.org/pdf/2405.18047"
